$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kHAP_diff_means")

# Header row
$ws.Cells.Item(1,1).Value = "kHAP"
$ws.Cells.Item(1,2).Value = "emmean"
$ws.Cells.Item(1,3).Value = "SE"
$ws.Cells.Item(1,4).Value = "df"
$ws.Cells.Item(1,5).Value = "lower.CL"
$ws.Cells.Item(1,6).Value = "upper.CL"

# Data rows
$ws.Cells.Item(2,1).Value = 100.0
$ws.Cells.Item(2,2).Value = 0.4501473
$ws.Cells.Item(2,3).Value = 0.01319982
$ws.Cells.Item(2,4).Value = 900.0
$ws.Cells.Item(2,5).Value = 0.4242413
$ws.Cells.Item(2,6).Value = 0.4760533

$ws.Cells.Item(3,1).Value = 250.0
$ws.Cells.Item(3,2).Value = 0.4318418
$ws.Cells.Item(3,3).Value = 0.01319982
$ws.Cells.Item(3,4).Value = 900.0
$ws.Cells.Item(3,5).Value = 0.4059358
$ws.Cells.Item(3,6).Value = 0.4577478

$ws.Cells.Item(4,1).Value = 500.0
$ws.Cells.Item(4,2).Value = 0.426946
$ws.Cells.Item(4,3).Value = 0.01319982
$ws.Cells.Item(4,4).Value = 900.0
$ws.Cells.Item(4,5).Value = 0.40104
$ws.Cells.Item(4,6).Value = 0.452852

$ws.Cells.Item(5,1).Value = 1000.0
$ws.Cells.Item(5,2).Value = 0.4177561
$ws.Cells.Item(5,3).Value = 0.01319982
$ws.Cells.Item(5,4).Value = 900.0
$ws.Cells.Item(5,5).Value = 0.3918501
$ws.Cells.Item(5,6).Value = 0.4436621

$ws.Cells.Item(6,1).Value = 2500.0
$ws.Cells.Item(6,2).Value = 0.3901631
$ws.Cells.Item(6,3).Value = 0.01319982
$ws.Cells.Item(6,4).Value = 900.0
$ws.Cells.Item(6,5).Value = 0.3642571
$ws.Cells.Item(6,6).Value = 0.4160691

$ws.Cells.Item(7,1).Value = 5000.0
$ws.Cells.Item(7,2).Value = 0.2984598
$ws.Cells.Item(7,3).Value = 0.01319982
$ws.Cells.Item(7,4).Value = 900.0
$ws.Cells.Item(7,5).Value = 0.2725538
$ws.Cells.Item(7,6).Value = 0.3243658

$ws.Cells.Item(8,1).Value = 10000.0
$ws.Cells.Item(8,2).Value = 0.2102065
$ws.Cells.Item(8,3).Value = 0.01319982
$ws.Cells.Item(8,4).Value = 900.0
$ws.Cells.Item(8,5).Value = 0.1843005
$ws.Cells.Item(8,6).Value = 0.2361125

$ws.Cells.Item(9,1).Value = 20000.0
$ws.Cells.Item(9,2).Value = 0.1813397
$ws.Cells.Item(9,3).Value = 0.01319982
$ws.Cells.Item(9,4).Value = 900.0
$ws.Cells.Item(9,5).Value = 0.1554337
$ws.Cells.Item(9,6).Value = 0.2072457

$ws.Cells.Item(10,1).Value = 30000.0
$ws.Cells.Item(10,2).Value = 0.1767917
$ws.Cells.Item(10,3).Value = 0.01319982
$ws.Cells.Item(10,4).Value = 900.0
$ws.Cells.Item(10,5).Value = 0.1508857
$ws.Cells.Item(10,6).Value = 0.2026977

# Centre-align the whole populated range, matching the style used elsewhere in the workbook
$ws.Range("A1:F10").HorizontalAlignment = -4108

# Make this sheet the active tab with A1:F10 selected (moves tabSelected off kHAP_diff_cont)
$ws.Activate()
$ws.Range("A1:F10").Select()
